# Apply updated crypto price/volume data to Sheet1
# Helper: sets a cell's value as text, forcing a quote-prefix when the
# string would otherwise be auto-converted to a number by Excel.
function Set-CellText {
    param($Ws, [string]$Addr, [string]$Val)
    if ($Val -match '^-?\d+(\.\d+)?$') {
        $Ws.Range($Addr).Value = "'" + $Val
    } else {
        $Ws.Range($Addr).Value = $Val
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "41.112.79"
Set-CellText $ws "E2" "  -1.80%  "

Set-CellText $ws "D3" "2.175.54"

Set-CellText $ws "D4" "0.999"
Set-CellText $ws "E4" "  -0.17%  "

Set-CellText $ws "D5" "248.44"
Set-CellText $ws "E5" "  -1.03%  "

Set-CellText $ws "D6" "0.614"
Set-CellText $ws "E6" "  -2.86%  "

Set-CellText $ws "D7" "65.38"
Set-CellText $ws "E7" "  -8.69%  "

Set-CellText $ws "E8" "  -0.09%  "

Set-CellText $ws "D9" "0.564"
Set-CellText $ws "E9" "  -4.87%  "

Set-CellText $ws "D10" "59.34"
Set-CellText $ws "E10" "  +1.78%  "

Set-CellText $ws "D11" "0.0926"
Set-CellText $ws "E11" "  -4.66%  "

Set-CellText $ws "D12" "35.49"
Set-CellText $ws "E12" "  -14.02%  "

Set-CellText $ws "E13" "  -1.42%  "

Set-CellText $ws "D14" "6.87"
Set-CellText $ws "E14" "  -5.08%  "

Set-CellText $ws "D15" "2.501.00"
Set-CellText $ws "E15" "  -2.25%  "

Set-CellText $ws "D16" "14.27"
Set-CellText $ws "E16" "  -4.94%  "

Set-CellText $ws "D17" "0.845"
Set-CellText $ws "E17" "  -2.51%  "

Set-CellText $ws "D18" "2.181.52"
Set-CellText $ws "E18" "  -2.16%  "

Set-CellText $ws "D19" "41.051.16"
Set-CellText $ws "E19" "  -1.84%  "

Set-CellText $ws "D20" "0.0₃0938"
Set-CellText $ws "E20" "  -3.87%  "

Set-CellText $ws "D21" "6.08"
Set-CellText $ws "E21" "  -2.43%  "

Set-CellText $ws "D22" "71.44"
Set-CellText $ws "E22" "  -2.17%  "

Set-CellText $ws "D23" "229.47"

Set-CellText $ws "E24" "  -4.36%  "

Set-CellText $ws "D25" "3.82"
Set-CellText $ws "E25" "  -9.62%  "

Set-CellText $ws "D27" "11.23"
Set-CellText $ws "E27" "  +4.22%  "

Set-CellText $ws "B28" "Toncoin"
Set-CellText $ws "C28" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-CellText $ws "D28" "2.45"
Set-CellText $ws "E28" "  +11.45%  "

Set-CellText $ws "B29" "PancakeSwap"
Set-CellText $ws "C29" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-CellText $ws "D29" "2.41"
Set-CellText $ws "E29" "  -5.87%  "

Set-CellText $ws "B30" "LEO"
Set-CellText $ws "C30" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-CellText $ws "D30" "3.72"
Set-CellText $ws "E30" "  -5.58%  "

Set-CellText $ws "D31" "167.98"
Set-CellText $ws "E31" "  -1.84%  "

Set-CellText $ws "E32" "  -3.20%  "

Set-CellText $ws "E33" "  -0.28%  "

Set-CellText $ws "D34" "5.67"
Set-CellText $ws "E34" "  +0.75%  "

Set-CellText $ws "E35" "  +1.40%  "

Set-CellText $ws "E36" "  -3.89%  "

Set-CellText $ws "E37" "  -4.42%  "

Set-CellText $ws "D38" "3.98"
Set-CellText $ws "E38" "  +1.29%  "

Set-CellText $ws "D39" "24.22"
Set-CellText $ws "E39" "  -6.91%  "

Set-CellText $ws "D40" "0.0305"
Set-CellText $ws "E40" "  +1.16%  "

Set-CellText $ws "E41" "  -5.24%  "

Set-CellText $ws "D42" "5.45"
Set-CellText $ws "E42" "  -8.81%  "

Set-CellText $ws "D43" "4.91"
Set-CellText $ws "E43" "  +0.03%  "

Set-CellText $ws "D44" "60.24"
Set-CellText $ws "E44" "  -11.32%  "

Set-CellText $ws "D45" "11.15"
Set-CellText $ws "E45" "  -5.71%  "

Set-CellText $ws "E46" "  -8.80%  "

Set-CellText $ws "D47" "8.49"
Set-CellText $ws "E47" "  -3.83%  "

Set-CellText $ws "B48" "Cronos"
Set-CellText $ws "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText $ws "D48" "0.0989"
Set-CellText $ws "E48" "  -3.34%  "

Set-CellText $ws "B49" "BinanceUSD"
Set-CellText $ws "C49" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-CellText $ws "D49" "1.00"
Set-CellText $ws "E49" "  -0.20%  "

Set-CellText $ws "D50" "1.15"
Set-CellText $ws "E50" "  -0.86%  "

Set-CellText $ws "D51" "1.15"
Set-CellText $ws "E51" "  -4.29%  "
